$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet (row 7: 73bf1b77-... handback report) ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "73bf1b77-69db-41ee-a069-6f357afba4d0.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a238187ca298da0acdbbc8b12f8b6f7bcb1d8f63/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md", "", "", "73bf1b77-69db-41ee-a069-6f357afba4d0.md")

$wsZh.Range("J7").Value = "73bf1b77-69db-41ee-a069-6f357afba4d0.7c4aed128300dc56cdc47e680606fa0187cb8b86.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-18 00:51:55"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6742e01d929cd3a0f849e82522112c42545b8078/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a238187ca298da0acdbbc8b12f8b6f7bcb1d8f63/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md."

# ---- de-de sheet (row 7: 73bf1b77-... handback report) ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "73bf1b77-69db-41ee-a069-6f357afba4d0.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a238187ca298da0acdbbc8b12f8b6f7bcb1d8f63/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md", "", "", "73bf1b77-69db-41ee-a069-6f357afba4d0.md")

$wsDe.Range("J7").Value = "73bf1b77-69db-41ee-a069-6f357afba4d0.7c4aed128300dc56cdc47e680606fa0187cb8b86.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-18 00:52:06"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6742e01d929cd3a0f849e82522112c42545b8078/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a238187ca298da0acdbbc8b12f8b6f7bcb1d8f63/e2e/73bf1b77-69db-41ee-a069-6f357afba4d0.md."

Write-Host "Generate Report for Handback: done"
